# edit.ps1 -- CryCompanywiseStockReport_1.xlsx data correction
#
# Several product rows had their per-batch figures (batch code B, rate E,
# quantity F, and amount G = D*F) swapped between adjacent rows of the same
# product. This script restores the correct figures for each affected row,
# one quantity correction (row 815), and recomputes the dependent Sub Total /
# Grand Total cells (which are stored as static numbers, not formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CHO-Medimix Sandal with Eladi oils for glowing skin and natural protection Soap-125gms  (rows 149, 150)
$ws.Range("B149").Value = 48654; $ws.Range("E149").Value = 38.26; $ws.Range("F149").Value = -1; $ws.Range("G149").Value = -32.02
$ws.Range("B150").Value = 63902; $ws.Range("E150").Value = 34.04; $ws.Range("F150").Value = 2; $ws.Range("G150").Value = 64.04000000000001
# DAB-Real Activ Coconut Water Tetra 1000ml  (rows 183, 184)
$ws.Range("B183").Value = 57552; $ws.Range("E183").Value = 136.86; $ws.Range("F183").Value = -5; $ws.Range("G183").Value = -603.45
$ws.Range("B184").Value = 64329; $ws.Range("E184").Value = 128.32; $ws.Range("F184").Value = 4; $ws.Range("G184").Value = 482.76
# HIM-GENTLE BABY SOAP 75G  (rows 279, 280)
$ws.Range("B279").Value = 48706; $ws.Range("E279").Value = 39.8; $ws.Range("F279").Value = -144; $ws.Range("G279").Value = -4795.2
$ws.Range("B280").Value = 64973; $ws.Range("E280").Value = 35.4; $ws.Range("F280").Value = 144; $ws.Range("G280").Value = 4795.2
# HUL-3Roses Dust [C] 500G Relaunch  (rows 313, 314)
$ws.Range("B313").Value = 57854; $ws.Range("F313").Value = 2; $ws.Range("G313").Value = 611.6799999999999
$ws.Range("B314").Value = 62997; $ws.Range("F314").Value = 0; $ws.Range("G314").Value = 0
# HUL-Bru Inst Poly 50g  (rows 316, 317, 318)
$ws.Range("B316").Value = 57077; $ws.Range("D316").Value = 93.08; $ws.Range("E316").Value = 111.2; $ws.Range("F316").Value = 1; $ws.Range("G316").Value = 93.08
$ws.Range("B317").Value = 61610; $ws.Range("D317").Value = 102.71; $ws.Range("E317").Value = 122.71; $ws.Range("F317").Value = -58; $ws.Range("G317").Value = -5957.18
$ws.Range("B318").Value = 63565; $ws.Range("E318").Value = 109.19; $ws.Range("F318").Value = 60; $ws.Range("G318").Value = 6162.6
# HUL-Kissan Pineapple Jam 500G  (rows 351, 352)
$ws.Range("B351").Value = 63531; $ws.Range("E351").Value = 152.53; $ws.Range("F351").Value = 80; $ws.Range("G351").Value = 11478.4
$ws.Range("B352").Value = 57802; $ws.Range("E352").Value = 162.71; $ws.Range("F352").Value = -79; $ws.Range("G352").Value = -11334.92
# HUL-Liril Soap 125 G  (rows 372, 373)
$ws.Range("B372").Value = 57885; $ws.Range("E372").Value = 62.28; $ws.Range("F372").Value = 4; $ws.Range("G372").Value = 208.52
$ws.Range("B373").Value = 63652; $ws.Range("E373").Value = 55.42; $ws.Range("F373").Value = 162; $ws.Range("G373").Value = 8445.059999999999
# HUL-lux advanced eventoned glow 4x100  (rows 375, 376)
$ws.Range("B375").Value = 63563; $ws.Range("E375").Value = 119.04; $ws.Range("F375").Value = 2; $ws.Range("G375").Value = 223.92
$ws.Range("B376").Value = 61605; $ws.Range("E376").Value = 133.78; $ws.Range("F376").Value = -13; $ws.Range("G376").Value = -1455.48
# HUL-Sfxl Ew Bale 500G  (rows 400, 401)
$ws.Range("B400").Value = 57835; $ws.Range("F400").Value = 1; $ws.Range("G400").Value = 59.13
$ws.Range("B401").Value = 62933; $ws.Range("F401").Value = 115; $ws.Range("G401").Value = 6799.95
# HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp  (rows 419, 420)
$ws.Range("B419").Value = 63007; $ws.Range("F419").Value = 834; $ws.Range("G419").Value = 142889.22
$ws.Range("B420").Value = 57856; $ws.Range("F420").Value = 2; $ws.Range("G420").Value = 342.66
# JLM-MBD Shiny Toothbrush Safari  (rows 457, 458)
$ws.Range("B457").Value = 63681; $ws.Range("E457").Value = 23.84; $ws.Range("F457").Value = 1; $ws.Range("G457").Value = 22.42
$ws.Range("B458").Value = 31930; $ws.Range("E458").Value = 26.8; $ws.Range("F458").Value = -62; $ws.Range("G458").Value = -1390.04
# CRE-Bourbon 100gm  (rows 579, 580)
$ws.Range("B579").Value = 65069; $ws.Range("E579").Value = 14.3; $ws.Range("F579").Value = 47; $ws.Range("G579").Value = 632.15
$ws.Range("B580").Value = 53757; $ws.Range("E580").Value = 16.08; $ws.Range("F580").Value = -159; $ws.Range("G580").Value = -2138.55
# CRE-Butter cookies 100gm  (rows 581, 582)
$ws.Range("B581").Value = 53602; $ws.Range("E581").Value = 15.69; $ws.Range("F581").Value = -231; $ws.Range("G581").Value = -3037.65
$ws.Range("B582").Value = 65068; $ws.Range("E582").Value = 13.97; $ws.Range("F582").Value = 181; $ws.Range("G582").Value = 2380.15
# CRE-Cremica Golden Bytes Rich Butter 200Gm  (rows 590, 591)
$ws.Range("B590").Value = 64922; $ws.Range("E590").Value = 20.98; $ws.Range("F590").Value = 165; $ws.Range("G590").Value = 3255.45
$ws.Range("B591").Value = 45706; $ws.Range("E591").Value = 23.58; $ws.Range("F591").Value = -202; $ws.Range("G591").Value = -3985.46
# CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm  (rows 593, 594)
$ws.Range("B593").Value = 45718; $ws.Range("E593").Value = 19.38; $ws.Range("F593").Value = -294; $ws.Range("G593").Value = -4768.68
$ws.Range("B594").Value = 64927; $ws.Range("E594").Value = 17.26; $ws.Range("F594").Value = 265; $ws.Range("G594").Value = 4298.3
# CRE-Cremica Oatmeal Digestive 112.5 Gm  (rows 599, 600)
$ws.Range("B599").Value = 45709; $ws.Range("E599").Value = 15.69; $ws.Range("F599").Value = -300; $ws.Range("G599").Value = -3945
$ws.Range("B600").Value = 64925; $ws.Range("E600").Value = 13.97; $ws.Range("F600").Value = 269; $ws.Range("G600").Value = 3537.35
# CRE-Cremica Pista Almond Cookies (75 +25Gm)  (rows 601, 602)
$ws.Range("B601").Value = 64919; $ws.Range("E601").Value = 27.97; $ws.Range("F601").Value = 178; $ws.Range("G601").Value = 4681.4
$ws.Range("B602").Value = 45702; $ws.Range("E602").Value = 31.43; $ws.Range("F602").Value = -215; $ws.Range("G602").Value = -5654.5
# CRE-Kaju khz cookies 100 gm  (rows 604, 605)
$ws.Range("B604").Value = 65067; $ws.Range("E604").Value = 15.65; $ws.Range("F604").Value = 283; $ws.Range("G604").Value = 4168.59
$ws.Range("B605").Value = 53595; $ws.Range("E605").Value = 17.61; $ws.Range("F605").Value = -335; $ws.Range("G605").Value = -4934.55
# PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)  (rows 687, 688)
$ws.Range("B687").Value = 64810; $ws.Range("E687").Value = 291.22; $ws.Range("F687").Value = 7; $ws.Range("G687").Value = 1917.44
$ws.Range("B688").Value = 53319; $ws.Range("E688").Value = 310.64; $ws.Range("F688").Value = -6; $ws.Range("G688").Value = -1643.52
# Rasna 32 Glass Shikanji Nimbupani  (rows 709, 710)
$ws.Range("B709").Value = 60025; $ws.Range("E709").Value = 37.22; $ws.Range("F709").Value = -98; $ws.Range("G709").Value = -3217.34
$ws.Range("B710").Value = 64833; $ws.Range("E710").Value = 34.9; $ws.Range("F710").Value = 97; $ws.Range("G710").Value = 3184.51
# Rasna Insta Orange 500g  (rows 715, 716)
$ws.Range("B715").Value = 64836; $ws.Range("E715").Value = 104.71; $ws.Range("F715").Value = 6; $ws.Range("G715").Value = 591
$ws.Range("B716").Value = 60031; $ws.Range("E716").Value = 111.69; $ws.Range("F716").Value = -5; $ws.Range("G716").Value = -492.5
# Rasna Nagpur Orange (32 Glass)  (rows 720, 721)
$ws.Range("B720").Value = 60022; $ws.Range("E720").Value = 37.22; $ws.Range("F720").Value = -113; $ws.Range("G720").Value = -3709.79
$ws.Range("B721").Value = 64830; $ws.Range("E721").Value = 34.9; $ws.Range("F721").Value = 114; $ws.Range("G721").Value = 3742.62
# SRL-Arctic Water Heater 15 Ltr  (rows 815)
$ws.Range("F815").Value = 1; $ws.Range("G815").Value = 4781.06
# Sub Total for SURYA ROSHNI LTD section = SUM(G809:G830)
$ws.Range("B831").Value = 65731.87
# Sub Total of all Sub Totals
$ws.Range("B962").Value = 4279191.69
# Grand Total (mirrors the overall Sub Total above)
$ws.Range("B963").Value = 4279191.69
